# "Add references and more" — update the table header labels to LaTeX-style
# math notation and refresh the active-cell selection on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tab = [char]9

$ws.Range("D1").Value = $tab + '$\text{Ctr}_k$'
$ws.Range("C1").Value = $tab + '$f_k$'
$ws.Range("B1").Value = $tab + '$n_k$'

$ws.Activate()
$ws.Range("G10").Select()
